# Auto-generated edit script applying numeric cell updates described by the commit diff.
# The workbook contains no formulas in these ranges (values were refreshed by a scheduled
# market-data runner), so each change is a direct literal value write.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4658.6
$ws.Range("J40").Value = 6777.4443
$ws.Range("L40").Value = 6777.4443
$ws.Range("N40").Value = -7127.4443
$ws.Range("H113").Value = 3225.1
$ws.Range("I113").Value = 2893
$ws.Range("K113").Value = 2893
$ws.Range("M113").Value = 361
$ws.Range("H138").Value = 1898.8871
$ws.Range("I138").Value = 1345.0883
$ws.Range("J138").Value = 2571.3572
$ws.Range("K138").Value = 4035.2649
$ws.Range("L138").Value = 7714.071599999999
$ws.Range("M138").Value = 1104.7351
$ws.Range("N138").Value = -17994.0716

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2287.2307
$ws.Range("I2").Value = 999.9474
$ws.Range("K2").Value = 999.9474
$ws.Range("M2").Value = -886.9474
$ws.Range("H5").Value = 328.25
$ws.Range("I5").Value = 108
$ws.Range("J5").Value = 989
$ws.Range("K5").Value = 108
$ws.Range("L5").Value = 989
$ws.Range("M5").Value = 4
$ws.Range("N5").Value = -1213
$ws.Range("H32").Value = 2771.6667
$ws.Range("I32").Value = 2487.524
$ws.Range("K32").Value = 2487.524
$ws.Range("M32").Value = -2200.524
$ws.Range("H45").Value = 114499.39
$ws.Range("I45").Value = 155941.47
$ws.Range("K45").Value = 155941.47
$ws.Range("M45").Value = -155564.47
$ws.Range("H55").Value = 32132
$ws.Range("J55").Value = 32132
$ws.Range("L55").Value = 32132
$ws.Range("N55").Value = -32762
$ws.Range("H63").Value = 4891.4
$ws.Range("J63").Value = 8282
$ws.Range("L63").Value = 8282
$ws.Range("N63").Value = -9654
$ws.Range("H66").Value = 4891.4
$ws.Range("J66").Value = 8282
$ws.Range("L66").Value = 41410
$ws.Range("N66").Value = -48274
$ws.Range("H97").Value = 621.0769
$ws.Range("J97").Value = 400
$ws.Range("L97").Value = 400
$ws.Range("N97").Value = -1392
$ws.Range("H116").Value = 2287.2307
$ws.Range("I116").Value = 999.9474
$ws.Range("K116").Value = 999.9474
$ws.Range("M116").Value = 1294.0526
$ws.Range("H132").Value = 7448.1724
$ws.Range("I132").Value = 6360.875
$ws.Range("K132").Value = 19082.625
$ws.Range("M132").Value = -16552.625

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2287.2307
$ws.Range("I3").Value = 999.9474
$ws.Range("K3").Value = 999.9474
$ws.Range("M3").Value = -885.9474
$ws.Range("H4").Value = 328.25
$ws.Range("I4").Value = 108
$ws.Range("J4").Value = 989
$ws.Range("K4").Value = 108
$ws.Range("L4").Value = 989
$ws.Range("M4").Value = 7
$ws.Range("N4").Value = -1219
$ws.Range("H35").Value = 42561.5
$ws.Range("J35").Value = 45913.8
$ws.Range("L35").Value = 45913.8
$ws.Range("N35").Value = -46533.8
$ws.Range("H82").Value = 42000
$ws.Range("J82").Value = 42000
$ws.Range("L82").Value = 42000
$ws.Range("N82").Value = -42766
$ws.Range("H85").Value = 42000
$ws.Range("J85").Value = 42000
$ws.Range("L85").Value = 42000
$ws.Range("N85").Value = -44652
$ws.Range("H86").Value = 1368.381
$ws.Range("I86").Value = 1368.381
$ws.Range("K86").Value = 1368.381
$ws.Range("M86").Value = -245.3810000000001
$ws.Range("H89").Value = 1368.381
$ws.Range("I89").Value = 1368.381
$ws.Range("K89").Value = 6841.905000000001
$ws.Range("M89").Value = -1225.905000000001
$ws.Range("H94").Value = 1655.8148
$ws.Range("I94").Value = 1519.381
$ws.Range("K94").Value = 1519.381
$ws.Range("M94").Value = -1068.381
$ws.Range("H134").Value = 3892.2144
$ws.Range("I134").Value = 3614.7437
$ws.Range("K134").Value = 10844.2311
$ws.Range("M134").Value = -8309.231100000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 23249.5
$ws.Range("J41").Value = 22300
$ws.Range("L41").Value = 22300
$ws.Range("N41").Value = -23156
$ws.Range("H51").Value = 29200
$ws.Range("J51").Value = 29200
$ws.Range("L51").Value = 29200
$ws.Range("N51").Value = -30672
$ws.Range("H60").Value = 27300
$ws.Range("J60").Value = 27300
$ws.Range("L60").Value = 27300
$ws.Range("N60").Value = -28322
$ws.Range("H61").Value = 29200
$ws.Range("J61").Value = 29200
$ws.Range("L61").Value = 29200
$ws.Range("N61").Value = -29896
$ws.Range("H68").Value = 40518.25
$ws.Range("J68").Value = 42357.668
$ws.Range("L68").Value = 42357.668
$ws.Range("N68").Value = -43855.668
$ws.Range("H71").Value = 40518.25
$ws.Range("J71").Value = 42357.668
$ws.Range("L71").Value = 127073.004
$ws.Range("N71").Value = -134561.004
$ws.Range("H74").Value = 38417.57
$ws.Range("J74").Value = 41439.668
$ws.Range("L74").Value = 41439.668
$ws.Range("N74").Value = -43187.668
$ws.Range("H77").Value = 38417.57
$ws.Range("J77").Value = 41439.668
$ws.Range("L77").Value = 124319.004
$ws.Range("N77").Value = -133055.004
$ws.Range("H134").Value = 8130.1763
$ws.Range("I134").Value = 7183.25
$ws.Range("J134").Value = 10402.8
$ws.Range("K134").Value = 21549.75
$ws.Range("L134").Value = 31208.4
$ws.Range("M134").Value = -19014.75
$ws.Range("N134").Value = -36278.39999999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2222
$ws.Range("I68").Value = 2849
$ws.Range("J68").Value = 1124.75
$ws.Range("K68").Value = 8547
$ws.Range("L68").Value = 3374.25
$ws.Range("M68").Value = -7736
$ws.Range("N68").Value = -4996.25
$ws.Range("H71").Value = 2222
$ws.Range("I71").Value = 2849
$ws.Range("J71").Value = 1124.75
$ws.Range("K71").Value = 25641
$ws.Range("L71").Value = 10122.75
$ws.Range("M71").Value = -21585
$ws.Range("N71").Value = -18234.75
$ws.Range("H94").Value = 900
$ws.Range("J94").Value = 900
$ws.Range("L94").Value = 2700
$ws.Range("N94").Value = -4052

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 51657.5
$ws.Range("J46").Value = 51657.5
$ws.Range("L46").Value = 51657.5
$ws.Range("N46").Value = -51969.5
$ws.Range("H136").Value = 46637.727
$ws.Range("J136").Value = 46637.727
$ws.Range("L136").Value = 139913.181
$ws.Range("N136").Value = -145013.181

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3462.375
$ws.Range("I7").Value = 3340
$ws.Range("K7").Value = 3340
$ws.Range("M7").Value = -3228
$ws.Range("H40").Value = 2886.625
$ws.Range("I40").Value = 2870.4285
$ws.Range("K40").Value = 2870.4285
$ws.Range("M40").Value = -2734.4285
$ws.Range("H122").Value = 3128.182
$ws.Range("I122").Value = 2849.8462
$ws.Range("K122").Value = 8549.5386
$ws.Range("M122").Value = -6099.5386
$ws.Range("H126").Value = 3462.375
$ws.Range("I126").Value = 3340
$ws.Range("K126").Value = 10020
$ws.Range("M126").Value = -7550
$ws.Range("H132").Value = 8044.614
$ws.Range("I132").Value = 7821.3687
$ws.Range("K132").Value = 23464.1061
$ws.Range("M132").Value = -20934.1061
$ws.Range("H136").Value = 4485.9614
$ws.Range("J136").Value = 7663.6665
$ws.Range("L136").Value = 22990.9995
$ws.Range("N136").Value = -28090.9995

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 29999
$ws.Range("J47").Value = 29999
$ws.Range("L47").Value = 29999
$ws.Range("N47").Value = -31143
$ws.Range("H51").Value = 31999
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H132").Value = 7721.8887
$ws.Range("I132").Value = 7249.5
$ws.Range("J132").Value = 8666.666999999999
$ws.Range("K132").Value = 21748.5
$ws.Range("L132").Value = 26000.001
$ws.Range("M132").Value = -19218.5
$ws.Range("N132").Value = -31060.001

